$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 343.63635
$ws.Range("I2").Value = 318.57144
$ws.Range("J2").Value = 387.5
$ws.Range("K2").Value = 318.57144
$ws.Range("L2").Value = 387.5
$ws.Range("M2").Value = -205.57144
$ws.Range("N2").Value = -613.5

# Row 31
$ws.Range("H31").Value = 22646.666
$ws.Range("I31").Value = 22646.666
$ws.Range("K31").Value = 67939.99800000001
$ws.Range("M31").Value = -67709.99800000001

# Row 42
$ws.Range("H42").Value = 576
$ws.Range("I42").Value = 1000
$ws.Range("J42").Value = 470
$ws.Range("K42").Value = 3000
$ws.Range("L42").Value = 1410
$ws.Range("M42").Value = -2770
$ws.Range("N42").Value = -1870

# Row 127
$ws.Range("H127").Value = 1136.5294
$ws.Range("I127").Value = 615.125
$ws.Range("J127").Value = 1600
$ws.Range("K127").Value = 1845.375
$ws.Range("L127").Value = 4800
$ws.Range("M127").Value = 3114.625
$ws.Range("N127").Value = -14720

# Row 137
$ws.Range("H137").Value = 1155.5883
$ws.Range("I137").Value = 1142.0256
$ws.Range("K137").Value = 3426.0768
$ws.Range("M137").Value = -876.0767999999998

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 110
$ws.Range("I5").Value = 110
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 110
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 2
$ws.Range("N5").ClearContents()

# Row 52
$ws.Range("H52").Value = 44350
$ws.Range("J52").Value = 44350
$ws.Range("L52").Value = 44350
$ws.Range("N52").Value = -44986

# Row 61
$ws.Range("H61").Value = 2528.3809
$ws.Range("I61").Value = 2256
$ws.Range("J61").Value = 3400
$ws.Range("K61").Value = 2256
$ws.Range("L61").Value = 3400
$ws.Range("M61").Value = -2044
$ws.Range("N61").Value = -3824

# Row 74
$ws.Range("H74").Value = 885.3333
$ws.Range("I74").Value = 716.5806
$ws.Range("J74").Value = 1539.25
$ws.Range("K74").Value = 716.5806
$ws.Range("L74").Value = 1539.25
$ws.Range("M74").Value = 157.4194
$ws.Range("N74").Value = -3287.25

# Row 77
$ws.Range("H77").Value = 885.3333
$ws.Range("I77").Value = 716.5806
$ws.Range("J77").Value = 1539.25
$ws.Range("K77").Value = 3582.903
$ws.Range("L77").Value = 7696.25
$ws.Range("M77").Value = 785.0969999999998
$ws.Range("N77").Value = -16432.25

# Row 132
$ws.Range("H132").Value = 3821.9038
$ws.Range("I132").Value = 4500.086
$ws.Range("J132").Value = 2425.647
$ws.Range("K132").Value = 13500.258
$ws.Range("L132").Value = 7276.941
$ws.Range("M132").Value = -10970.258
$ws.Range("N132").Value = -12336.941

# Row 134
$ws.Range("H134").Value = 44997.5
$ws.Range("J134").Value = 44997.5
$ws.Range("L134").Value = 44997.5
$ws.Range("N134").Value = -55137.5

# Row 136
$ws.Range("H136").Value = 2528.3809
$ws.Range("I136").Value = 2256
$ws.Range("J136").Value = 3400
$ws.Range("K136").Value = 6768
$ws.Range("L136").Value = 10200
$ws.Range("M136").Value = -4218
$ws.Range("N136").Value = -15300

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 110
$ws.Range("I4").Value = 110
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 110
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 5
$ws.Range("N4").ClearContents()

# Row 22
$ws.Range("H22").Value = 4909.1816
$ws.Range("I22").Value = 4909.1816
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 4909.1816
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -4736.1816
$ws.Range("N22").ClearContents()

# Row 134
$ws.Range("H134").Value = 2774.647
$ws.Range("I134").Value = 2640.65
$ws.Range("J134").Value = 2966.0715
$ws.Range("K134").Value = 7921.950000000001
$ws.Range("L134").Value = 8898.2145
$ws.Range("M134").Value = -5386.950000000001
$ws.Range("N134").Value = -13968.2145

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4750
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

# Row 34
$ws.Range("H34").Value = 4750
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

# Row 58
$ws.Range("H58").Value = 757595.1
$ws.Range("I58").Value = 1159007.6
$ws.Range("J58").Value = 1995.0588
$ws.Range("K58").Value = 1159007.6
$ws.Range("L58").Value = 1995.0588
$ws.Range("M58").Value = -1158804.6
$ws.Range("N58").Value = -2401.0588

# Row 132
$ws.Range("H132").Value = 484382.1
$ws.Range("I132").Value = 645089.5
$ws.Range("J132").Value = 2260
$ws.Range("K132").Value = 1935268.5
$ws.Range("L132").Value = 6780
$ws.Range("M132").Value = -1932738.5
$ws.Range("N132").Value = -11840

# Row 134
$ws.Range("H134").Value = 3553.25
$ws.Range("I134").Value = 2953
$ws.Range("J134").Value = 3853.375
$ws.Range("K134").Value = 8859
$ws.Range("L134").Value = 11560.125
$ws.Range("M134").Value = -6324
$ws.Range("N134").Value = -16630.125

# Row 136
$ws.Range("H136").Value = 757595.1
$ws.Range("I136").Value = 1159007.6
$ws.Range("J136").Value = 1995.0588
$ws.Range("K136").Value = 3477022.8
$ws.Range("L136").Value = 5985.1764
$ws.Range("M136").Value = -3474472.8
$ws.Range("N136").Value = -11085.1764

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 218.7
$ws.Range("J23").Value = 244.625
$ws.Range("L23").Value = 733.875
$ws.Range("N23").Value = -1203.875

# Row 76
$ws.Range("H76").Value = 2802.6
$ws.Range("I76").Value = 2013
$ws.Range("K76").Value = 6039
$ws.Range("M76").Value = -5656

# Row 79
$ws.Range("H79").Value = 2802.6
$ws.Range("I79").Value = 2013
$ws.Range("K79").Value = 6039
$ws.Range("M79").Value = -4713

# Row 115
$ws.Range("H115").Value = 1914
$ws.Range("I115").Value = 1914
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 5742
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -4567
$ws.Range("N115").ClearContents()

# Row 117
$ws.Range("H117").Value = 770.4
$ws.Range("J117").Value = 796.3333
$ws.Range("L117").Value = 2388.9999
$ws.Range("N117").Value = -9272.999899999999

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2369.0605
$ws.Range("I132").Value = 1687.3529
$ws.Range("J132").Value = 3093.375
$ws.Range("K132").Value = 5062.0587
$ws.Range("L132").Value = 9280.125
$ws.Range("M132").Value = -2532.0587
$ws.Range("N132").Value = -14340.125

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 763.1
$ws.Range("I22").Value = 808
$ws.Range("J22").Value = 726.36365
$ws.Range("K22").Value = 808
$ws.Range("L22").Value = 726.36365
$ws.Range("M22").Value = -513
$ws.Range("N22").Value = -1316.36365

# Row 27
$ws.Range("H27").Value = 763.1
$ws.Range("I27").Value = 808
$ws.Range("J27").Value = 726.36365
$ws.Range("K27").Value = 808
$ws.Range("L27").Value = 726.36365
$ws.Range("M27").Value = -701
$ws.Range("N27").Value = -940.36365

# Row 132
$ws.Range("H132").Value = 2519.2444
$ws.Range("I132").Value = 2148.121
$ws.Range("J132").Value = 3539.8333
$ws.Range("K132").Value = 6444.363
$ws.Range("L132").Value = 10619.4999
$ws.Range("M132").Value = -3914.363
$ws.Range("N132").Value = -15679.4999

# Row 136
$ws.Range("H136").Value = 27299482
$ws.Range("I136").Value = 35715784
$ws.Range("J136").Value = 1115434.5
$ws.Range("K136").Value = 107147352
$ws.Range("L136").Value = 3346303.5
$ws.Range("M136").Value = -107144802
$ws.Range("N136").Value = -3351403.5

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1603.551
$ws.Range("I132").Value = 1031.8064
$ws.Range("J132").Value = 2588.2222
$ws.Range("K132").Value = 3095.4192
$ws.Range("L132").Value = 7764.6666
$ws.Range("M132").Value = -565.4191999999998
$ws.Range("N132").Value = -12824.6666

# Row 136
$ws.Range("H136").Value = 1729.1428
$ws.Range("I136").Value = 1564.7894
$ws.Range("J136").Value = 2076.111
$ws.Range("K136").Value = 4694.3682
$ws.Range("L136").Value = 6228.333
$ws.Range("M136").Value = -2144.3682
$ws.Range("N136").Value = -11328.333
